# Apply cryptos price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each target cell to Text format first so Excel keeps these
# numeric-looking / percent strings as literal text (matching the
# original inlineStr cells) instead of re-parsing them as numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.571.50'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.975.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -5.19%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.33'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -5.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.23'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -8.32%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.986.03'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.31%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.14'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -6.90%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.497.99'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.07%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '61.640.01'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.65'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.979.71'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.11%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -5.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.16'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.05'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.05'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.51%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.48%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.73'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.471'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.094.90'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.57%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0941'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.55%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.83%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.45'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.98'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.86%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.71%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.55%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.52%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.53'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.411.48'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.11'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.671'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0593'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.997'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.06'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.71%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0954'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.78'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.55%  '
